$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": update Version value from "3.2.0-ballot" to "3.2.0" ---
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B3").Value = "3.2.0"

# --- Sheet "Include from ch-ehealth-codes": clear B6 ("Document search") ---
$wsInclude = $wb.Worksheets.Item("Include from ch-ehealth-codes")
$wsInclude.Range("B6").ClearContents()
